$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.060.33"
$ws.Range("E2").Value = "  +0.09%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.874.38"
$ws.Range("E3").Value = "  -1.28%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.57"
$ws.Range("E5").Value = "  +0.34%  "

$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5068"
$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3840"
$ws.Range("E8").Value = "  -1.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08957"
$ws.Range("E9").Value = "  -2.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.120"
$ws.Range("E10").Value = "  -1.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.61"
$ws.Range("E11").Value = "  -0.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.327"
$ws.Range("E12").Value = "  -0.37%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.71"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.869.92"
$ws.Range("E14").Value = "  -1.77%  "

$ws.Range("B15").Value = "BinanceUSD"
$ws.Range("C15").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.004"
$ws.Range("E15").Value = "  +0.17%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.183"
$ws.Range("E16").Value = "  -1.38%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001105"
$ws.Range("E17").Value = "  -1.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.03"
$ws.Range("E18").Value = "  -1.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06595"
$ws.Range("E19").Value = "  +0.40%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.13"
$ws.Range("E20").Value = "  +2.23%  "

$ws.Range("E21").Value = "  +0.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.105"
$ws.Range("E22").Value = "  -1.58%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.081.39"
$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.40"
$ws.Range("E24").Value = "  +0.53%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.280"
$ws.Range("E25").Value = "  -1.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.089.73"
$ws.Range("E26").Value = "  -1.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.534"
$ws.Range("E27").Value = "  -2.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.75"
$ws.Range("E28").Value = "  -0.69%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "157.08"
$ws.Range("E29").Value = "  -0.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.55"
$ws.Range("E30").Value = "  -0.28%  "

$ws.Range("E31").Value = "  -1.56%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.057"
$ws.Range("E32").Value = "  -2.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.603"
$ws.Range("E33").Value = "  +0.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.604"
$ws.Range("E34").Value = "  -0.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.639"
$ws.Range("E35").Value = "  +0.73%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06572"
$ws.Range("E36").Value = "  -1.33%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02425"
$ws.Range("E37").Value = "  +0.95%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2175"
$ws.Range("E38").Value = "  +0.11%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.206"
$ws.Range("E39").Value = "  -1.72%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.266"
$ws.Range("E40").Value = "  +1.42%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6378"
$ws.Range("E41").Value = "  +0.42%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.41"
$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.905"

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.25"
$ws.Range("E44").Value = "  -0.31%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6008"
$ws.Range("E45").Value = "  +0.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.671"
$ws.Range("E46").Value = "  -0.87%  "

$ws.Range("E47").Value = "  +0.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.237"
$ws.Range("E48").Value = "  +5.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.990"
$ws.Range("E49").Value = "  -0.63%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "121.26"
$ws.Range("E50").Value = "  -1.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.61"
$ws.Range("E51").Value = "  +1.84%  "
